# Weekly driver report update for 2025-04-21
# Updates counts/percentages in the "Bad Drivers" table and swaps the
# data rows for the two drivers that changed rank order (row 7 <-> row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 4: Intel(R) Wi-Fi 6 AX201 160MHz - 22.200.2.1
$ws.Range("C4").Value = 4678
$ws.Range("D4").Value = 81.2

# Row 5: Intel(R) Wi-Fi 6 AX201 160MHz - 22.0.1.5
$ws.Range("C5").Value = 290

# Row 6: Intel(R) Wi-Fi 6 AX201 160MHz - 23.90.0.2
$ws.Range("C6").Value = 663
$ws.Range("D6").Value = 97.2

# Row 7 now holds the Intel(R) Wi-Fi 6E AX211 160MHz - 22.200.2.1 data
$ws.Range("A7").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.200.2.1"
$ws.Range("B7").Value = 43
$ws.Range("C7").Value = 1798
$ws.Range("D7").Value = 98.2

# Row 8 now holds the Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2 data
$ws.Range("A8").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2"
$ws.Range("B8").Value = 11
$ws.Range("C8").Value = 437
$ws.Range("D8").Value = 98.2

# Row 9: Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.0.3
$ws.Range("B9").Value = 243
$ws.Range("C9").Value = 7147

# Row 10: Intel(R) Wi-Fi 6 AX201 160MHz - 22.170.2.1
$ws.Range("B10").Value = 509
$ws.Range("C10").Value = 18236

# Row 11: Intel(R) Wi-Fi 6E AX211 160MHz - 22.170.2.1
$ws.Range("B11").Value = 286
$ws.Range("C11").Value = 5672

# Row 12: Totals
$ws.Range("B12").Value = 1110
$ws.Range("C12").Value = 38948
